$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = '52.032.74'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = '  +5.10%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = '2.779.38'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = '  +5.34%  '
$ws.Range("E4").Value2 = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '115.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = '  +2.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '340.20'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = '  +4.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = '0.548'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value2 = '  +4.71%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = '0.576'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = '  +4.84%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = '41.83'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = '  +5.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = '0.0860'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = '  +5.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = '20.04'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = '  +0.22%  '
$ws.Range("E13").Value2 = '  +2.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '7.60'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = '3.218.47'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = '  +5.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = '2.776.75'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = '  +5.44%  '
$ws.Range("B17").Value2 = 'WrappedBTC'
$ws.Range("C17").Value2 = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = '51.891.68'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = '  +4.85%  '
$ws.Range("B18").Value2 = 'Polygon'
$ws.Range("C18").Value2 = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = '0.878'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = '  +2.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '3.19'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = '  +9.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '13.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = '  -0.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = '6.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = '  +4.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '0.0₃0976'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = '  +2.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = '276.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = '  +3.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '69.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = '  +1.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = '2.74'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = '  +7.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = '26.67'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = '  +2.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = '10.19'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = '  +0.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = '2.22'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = '  +1.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = '0.141'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = '  +2.72%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = '34.74'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = '  +0.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '50.20'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = '  +1.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = '5.70'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = '  +4.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '0.0824'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = '  +0.94%  '
$ws.Range("E35").Value2 = '  -0.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = '2.10'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = '  +3.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = '18.82'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = '  -1.63%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = '4.93'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = '  -0.38%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = '3.23'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = '  +4.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '0.0378'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = '  +11.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = '2.66'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = '  +25.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '2.34'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = '  +1.28%  '
$ws.Range("E43").Value2 = '  +3.53%  '
$ws.Range("B44").Value2 = 'Monero'
$ws.Range("C44").Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = '125.53'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = '  -3.06%  '
$ws.Range("B45").Value2 = 'EnergySwap'
$ws.Range("C45").Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = '23.04'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = '  +0.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = '2.078.27'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = '  +0.77%  '
$ws.Range("E47").Value2 = '  +0.15%  '
$ws.Range("E48").Value2 = '  +3.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '5.53'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = '  +5.90%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = '8.90'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = '  +0.20%  '
$ws.Range("B51").Value2 = 'MultiversX'
$ws.Range("C51").Value2 = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = '59.33'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = '  +1.13%  '
